$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.783.06"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "3.497.35"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.30"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.52"
$ws.Range("E6").Value = "  -7.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +3.81%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.26"
$ws.Range("E11").Value = "  -5.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000271"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.17"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "4.055.08"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "3.497.47"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.30"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "65.793.33"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.997"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.02"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  +4.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.89"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.11"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.76"
$ws.Range("E25").Value = "  +6.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.86"
$ws.Range("E26").Value = "  -5.52%  "
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.95"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.15"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "618.03"
$ws.Range("E30").Value = "  -10.79%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("E31").Value = "  -7.32%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.63"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.10"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.151"
$ws.Range("E35").Value = "  +11.74%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0790"
$ws.Range("E37").Value = "  -5.24%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.97"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").Value = "3.370.85"
$ws.Range("E39").Value = "  +10.76%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.378"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.36"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0415"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.49"
$ws.Range("E46").Value = "  -8.79%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.132"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.39"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.41"
$ws.Range("E50").Value = "  -9.70%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.81"
$ws.Range("E51").Value = "  +7.55%  "
